$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.817.50'
$ws.Range('E2').Value = '  -1.33%  '
$ws.Range('D3').Value = '2.397.41'
$ws.Range('E3').Value = '  -2.25%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '552.81'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.92%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '158.02'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.67%  '
$ws.Range('E7').Value = '  +0.11%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('E10').Value = '  -1.47%  '
$ws.Range('E11').Value = '  -2.37%  '
$ws.Range('E12').Value = '  -2.23%  '
$ws.Range('D13').Value = '67.749.81'
$ws.Range('E13').Value = '  -1.11%  '
$ws.Range('E14').Value = '  -0.34%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '22.75'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.32%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '10.26'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -4.53%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '329.06'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -3.43%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.78'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -3.96%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.74'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.21%  '
$ws.Range('E20').Value = '  +0.12%  '
$ws.Range('E21').Value = '  -3.72%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '65.56'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.95%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.61'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.01%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.04'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.76%  '
$ws.Range('D25').Value = '0.0₃0792'
$ws.Range('E25').Value = '  -3.02%  '
$ws.Range('E26').Value = '  -1.89%  '
$ws.Range('E27').Value = '  +0.15%  '
$ws.Range('B28').Value = 'Fetch.AI'
$ws.Range('C28').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.13'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.39%  '
$ws.Range('B29').Value = 'Bittensor'
$ws.Range('C29').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '415.46'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -4.68%  '
$ws.Range('E30').Value = '  -1.58%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '157.15'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.00%  '
$ws.Range('E32').Value = '  -0.21%  '
$ws.Range('E33').Value = '  -0.01%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '17.63'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.95%  '
$ws.Range('E35').Value = '  -3.82%  '
$ws.Range('E36').Value = '  -3.18%  '
$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.21'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -5.35%  '
$ws.Range('B38').Value = 'Stacks'
$ws.Range('C38').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.45'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.57%  '
$ws.Range('E39').Value = '  -4.65%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.28'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.82%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '128.20'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.50%  '
$ws.Range('E42').Value = '  -7.40%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0705'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.91%  '
$ws.Range('E44').Value = '  -1.65%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.551'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.41%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0912'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.79%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.10'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.77%  '
$ws.Range('E48').Value = '  -6.50%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '16.34'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.24%  '
$ws.Range('E50').Value = '  -6.31%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0425'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.48%  '
